# Updated cryptos list on Sat Oct 14 08:41:37 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns of the crypto table with new
# quotes, and swaps the FraxShare / PaxDollar rows (41 <-> 42) to reflect the
# new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.904.50'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.550.03'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '206.49'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '22.16'
$ws.Range('E8').Value = '  +2.87%  '
$ws.Range('E9').Value = '  -0.83%  '
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').Value = '1.770.52'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '1.550.85'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.73'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.518'
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('D16').Value = '26.906.12'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.61'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '217.48'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.05'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.22'
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '154.25'
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  +1.62%  '
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.22'
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').Value = '1.415.78'
$ws.Range('E33').Value = '  +3.01%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.06'
$ws.Range('E34').Value = '  +3.54%  '
$ws.Range('E35').Value = '  +2.17%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.964'
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.524'
$ws.Range('E39').Value = '  +0.68%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.807'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.71'
$ws.Range('E42').Value = '  +3.45%  '
$ws.Range('E43').Value = '  +3.95%  '
$ws.Range('E44').Value = '  +1.76%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '64.52'
$ws.Range('E45').Value = '  +1.19%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.75'
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').Value = '1.684.56'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '87.58'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('E49').Value = '  +6.20%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0517'
$ws.Range('E50').Value = '  +1.82%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0960'
$ws.Range('E51').Value = '  +0.34%  '
